$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 11766.777
$ws.Range("I107").Value = 17150.834
$ws.Range("J107").Value = 998.6667
$ws.Range("K107").Value = 17150.834
$ws.Range("L107").Value = 998.6667
$ws.Range("M107").Value = -15230.834
$ws.Range("N107").Value = -4838.6667

$ws.Range("H138").Value = 1300.8776
$ws.Range("J138").Value = 2190.7778
$ws.Range("L138").Value = 6572.3334
$ws.Range("N138").Value = -16852.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1452477.8
$ws.Range("I63").Value = 3705959.2
$ws.Range("J63").Value = 3811.0715
$ws.Range("K63").Value = 3705959.2
$ws.Range("L63").Value = 3811.0715
$ws.Range("M63").Value = -3705273.2
$ws.Range("N63").Value = -5183.0715

$ws.Range("H66").Value = 1452477.8
$ws.Range("I66").Value = 3705959.2
$ws.Range("J66").Value = 3811.0715
$ws.Range("K66").Value = 18529796
$ws.Range("L66").Value = 19055.3575
$ws.Range("M66").Value = -18526364
$ws.Range("N66").Value = -25919.3575

$ws.Range("H88").Value = 11212.529
$ws.Range("I88").Value = 2612.375
$ws.Range("J88").Value = 18857.111
$ws.Range("K88").Value = 2612.375
$ws.Range("L88").Value = 18857.111
$ws.Range("M88").Value = -2206.375
$ws.Range("N88").Value = -19669.111

$ws.Range("H91").Value = 11212.529
$ws.Range("I91").Value = 2612.375
$ws.Range("J91").Value = 18857.111
$ws.Range("K91").Value = 2612.375
$ws.Range("L91").Value = 18857.111
$ws.Range("M91").Value = -1208.375
$ws.Range("N91").Value = -21665.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2728
$ws.Range("I86").Value = 2160
$ws.Range("K86").Value = 2160
$ws.Range("M86").Value = -1037

$ws.Range("H89").Value = 2728
$ws.Range("I89").Value = 2160
$ws.Range("K89").Value = 10800
$ws.Range("M89").Value = -5184

$ws.Range("H99").Value = 883.94446
$ws.Range("I99").Value = 758.82355
$ws.Range("J99").Value = 3011
$ws.Range("K99").Value = 758.82355
$ws.Range("L99").Value = 3011
$ws.Range("M99").Value = 739.17645
$ws.Range("N99").Value = -6007

$ws.Range("H105").Value = 2092.3462
$ws.Range("I105").Value = 2042.381
$ws.Range("J105").Value = 2302.2
$ws.Range("K105").Value = 2042.381
$ws.Range("L105").Value = 2302.2
$ws.Range("M105").Value = -295.3810000000001
$ws.Range("N105").Value = -5796.2

$ws.Range("H134").Value = 10103790
$ws.Range("I134").Value = 12823083
$ws.Range("K134").Value = 38469249
$ws.Range("M134").Value = -38466714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 91500.91
$ws.Range("I16").Value = 111662.22
$ws.Range("K16").Value = 111662.22
$ws.Range("M16").Value = -111375.22

$ws.Range("H31").Value = 10876.655
$ws.Range("I31").Value = 1191.6945
$ws.Range("K31").Value = 1191.6945
$ws.Range("M31").Value = -896.6945000000001

$ws.Range("H34").Value = 10876.655
$ws.Range("I34").Value = 1191.6945
$ws.Range("K34").Value = 1191.6945
$ws.Range("M34").Value = -989.6945000000001

$ws.Range("H58").Value = 2864453.2
$ws.Range("I58").Value = 4329637.5
$ws.Range("J58").Value = 20272.234
$ws.Range("K58").Value = 4329637.5
$ws.Range("L58").Value = 20272.234
$ws.Range("M58").Value = -4329434.5
$ws.Range("N58").Value = -20678.234

$ws.Range("H86").Value = 2521.9524
$ws.Range("I86").Value = 2350.9167
$ws.Range("J86").Value = 2750
$ws.Range("K86").Value = 2350.9167
$ws.Range("L86").Value = 2750
$ws.Range("M86").Value = -1227.9167
$ws.Range("N86").Value = -4996

$ws.Range("H89").Value = 2521.9524
$ws.Range("I89").Value = 2350.9167
$ws.Range("J89").Value = 2750
$ws.Range("K89").Value = 11754.5835
$ws.Range("L89").Value = 13750
$ws.Range("M89").Value = -6138.583500000001
$ws.Range("N89").Value = -24982

$ws.Range("H94").Value = 2054.4814
$ws.Range("I94").Value = 1224.6666
$ws.Range("J94").Value = 2291.5715
$ws.Range("K94").Value = 1224.6666
$ws.Range("L94").Value = 2291.5715
$ws.Range("M94").Value = -773.6666
$ws.Range("N94").Value = -3193.5715

$ws.Range("H107").Value = 267.1111
$ws.Range("I107").Value = 211.6923
$ws.Range("J107").Value = 411.2
$ws.Range("K107").Value = 211.6923
$ws.Range("L107").Value = 411.2
$ws.Range("M107").Value = 1708.3077
$ws.Range("N107").Value = -4251.2

$ws.Range("H113").Value = 91500.91
$ws.Range("I113").Value = 111662.22
$ws.Range("K113").Value = 111662.22
$ws.Range("M113").Value = -109492.22

$ws.Range("H136").Value = 2864453.2
$ws.Range("I136").Value = 4329637.5
$ws.Range("J136").Value = 20272.234
$ws.Range("K136").Value = 12988912.5
$ws.Range("L136").Value = 60816.702
$ws.Range("M136").Value = -12986362.5
$ws.Range("N136").Value = -65916.702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 183.58333
$ws.Range("I33").Value = 132.6
$ws.Range("J33").Value = 220
$ws.Range("K33").Value = 795.5999999999999
$ws.Range("L33").Value = 1320
$ws.Range("M33").Value = -512.5999999999999
$ws.Range("N33").Value = -1886

$ws.Range("H112").Value = 2448
$ws.Range("I112").Value = 1840
$ws.Range("J112").Value = 2600
$ws.Range("K112").Value = 5520
$ws.Range("L112").Value = 7800
$ws.Range("M112").Value = -4412
$ws.Range("N112").Value = -10016

$ws.Range("H115").Value = 5000
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H132").Value = 1890.8182
$ws.Range("I132").Value = 799.8
$ws.Range("K132").Value = 7198.2
$ws.Range("M132").Value = -4668.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25895.936
$ws.Range("I70").Value = 34312.363
$ws.Range("J70").Value = 4531.154
$ws.Range("K70").Value = 34312.363
$ws.Range("L70").Value = 4531.154
$ws.Range("M70").Value = -34042.363
$ws.Range("N70").Value = -5071.154

$ws.Range("H73").Value = 25895.936
$ws.Range("I73").Value = 34312.363
$ws.Range("J73").Value = 4531.154
$ws.Range("K73").Value = 34312.363
$ws.Range("L73").Value = 4531.154
$ws.Range("M73").Value = -33376.363
$ws.Range("N73").Value = -6403.154

$ws.Range("H107").Value = 77303.84
$ws.Range("I107").Value = 125378.75
$ws.Range("J107").Value = 384
$ws.Range("K107").Value = 125378.75
$ws.Range("L107").Value = 384
$ws.Range("M107").Value = -123458.75
$ws.Range("N107").Value = -4224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 16461.25
$ws.Range("I16").Value = 6422
$ws.Range("J16").Value = 26500.5
$ws.Range("K16").Value = 6422
$ws.Range("L16").Value = 26500.5
$ws.Range("M16").Value = -6252
$ws.Range("N16").Value = -26840.5

$ws.Range("H46").Value = 5365.9546
$ws.Range("I46").Value = 903.1875
$ws.Range("J46").Value = 17266.666
$ws.Range("K46").Value = 903.1875
$ws.Range("L46").Value = 17266.666
$ws.Range("M46").Value = -715.1875
$ws.Range("N46").Value = -17642.666

$ws.Range("H61").Value = 1900.4166
$ws.Range("I61").Value = 1316.6666
$ws.Range("J61").Value = 2484.1667
$ws.Range("K61").Value = 1316.6666
$ws.Range("L61").Value = 2484.1667
$ws.Range("M61").Value = -1114.6666
$ws.Range("N61").Value = -2888.1667

$ws.Range("H93").Value = 46934.047
$ws.Range("I93").Value = 672.64703
$ws.Range("J93").Value = 204222.8
$ws.Range("K93").Value = 672.64703
$ws.Range("L93").Value = 204222.8
$ws.Range("M93").Value = 575.35297
$ws.Range("N93").Value = -206718.8

$ws.Range("H113").Value = 1900.4166
$ws.Range("I113").Value = 1316.6666
$ws.Range("J113").Value = 2484.1667
$ws.Range("K113").Value = 1316.6666
$ws.Range("L113").Value = 2484.1667
$ws.Range("M113").Value = 853.3334
$ws.Range("N113").Value = -6824.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 159.5
$ws.Range("I100").Value = 148.21428
$ws.Range("J100").Value = 199
$ws.Range("K100").Value = 296.42856
$ws.Range("L100").Value = 398
$ws.Range("M100").Value = 244.57144
$ws.Range("N100").Value = -1480
